$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header date C1: 2018-02-01 -> 2018-02-12
$ws.Range("C1").Value = 43143

# Row 7: Battery charge view in Hola App -> status Done -> Closed
$ws.Range("C7").Value = "Closed"

# Row 32: Support for multiple buddies -> status Ongoing -> Closed, clear notes
$ws.Range("C32").Value = "Closed"
$ws.Range("D32").Clear()

# Row 56: Install stronger speaker -> Ongoing -> Closed
$ws.Range("C56").Value = "Closed"

# Row 58: Use exponential back-off for EchoReq -> Open -> Rejected
$ws.Range("C58").Value = "Rejected"

# Row 61: Experiment with ADPCM 16MHz mode -> Ongoing -> Closed
$ws.Range("C61").Value = "Closed"

# Row 64: Sync app to latest SDK -> Open -> Closed
$ws.Range("C64").Value = "Closed"

# Row 72: Redo server side design diagram in draw.io -> Open -> Rejected, add note
$ws.Range("C72").Value = "Rejected"
$ws.Range("D72").Value = "Using yuml instead"

# New rows 73-79
$ws.Range("B73").Value = "Add button to enter setup mode"
$ws.Range("C73").Value = "Open"

$ws.Range("B74").Value = "Scalable server architecture"
$ws.Range("C74").Value = "Open"

$ws.Range("B75").Value = "Buddy config in app takes too long"
$ws.Range("C75").Value = "Open"

$ws.Range("B76").Value = "Intercom gets stuck in UDP can't send state"
$ws.Range("C76").Value = "Open"

$ws.Range("B77").Value = "Volume Control not working well"
$ws.Range("C77").Value = "Open"

$ws.Range("B78").Value = "Hissing, crackling, high pitch noise come from speakers in rest"
$ws.Range("C78").Value = "Open"

$ws.Range("B79").Value = "Intercom2 speaker does not work"
$ws.Range("C79").Value = "Open"

# Update view: select C1 (also clears the scrolled topLeftCell from before)
$ws.Activate()
$ws.Range("C1").Select()
